$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the existing hyperlinks before changing the cell text
$ws.Hyperlinks.Delete()

# Update names
$ws.Range("B2").Value = "nama satu"
$ws.Range("B3").Value = "nama dua"
$ws.Range("B4").Value = "nama tiga"

# Update emails (cells with hyperlinks)
$ws.Range("D2").Value = "satu@mail.com"
$ws.Range("D3").Value = "dua@mail.com"
$ws.Range("D4").Value = "tiga@mail.com"

# Recreate the hyperlinks with new targets/display text
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:satu@mail.com", "", "", "satu@mail.com")
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:dua@mail.com", "", "", "dua@mail.com")
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:tiga@mail.com", "", "", "tiga@mail.com")

# Update selection to D5
$ws.Range("D5").Select()
